$d = $word.ActiveDocument

# --- Edit 1: "Decision Making" paragraph -------------------------------
$old1 = "We aim to reach a consensus once everyone has had a chance to voice their opinions. If it doesn't work, we might think about getting the TA's opinion first, followed by a majority vote. We believe that consensus will be reached in the majority of cases."
$new1 = "We aim to reach a consensus once everyone has had a chance to voice their opinions. If it doesn't work, we might think about getting the TA's opinion first, so that some members may change their opinion in order to facilitate a consensus. If this also does not work, we follow with a majority vote. The minority that disagrees with the vote will be obliged to follow the majority decision nonetheless. This decision will not be harmful for the project, since we would have asked the TA" + [char]0x2019 + "s opinion first. "

$r1 = $d.Content
$r1.Find.ClearFormatting()
$found1 = $r1.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $r1.Text = $new1
}

# --- Edit 2: "Dealing with Conflicts" paragraph -------------------------
# The old text is split across three separate runs in the original
# document; Find can still locate it as one contiguous range spanning the
# run boundaries, and re-assigning .Text collapses it into a single run.
$old2 = "Similarly to decision-making, we discuss the matter within the group, and follow the same process as described before, namely, prioritise consensus and fall back on majority vote. In order to reach a consensus, it is essential to ask for everyone" + [char]0x2019 + "s opinion to understand what points each member raises and thus make a decision based thereon. We may ask our TA for assistance if appropriate."
$new2 = "Similarly to decision-making, we discuss the matter within the group, and prioritise consensus. The discussion should be respectful and must not involve any personal matters. If the conflict persists after a discussion, we request help from our TA, so that the correct measures can be taken, which may include a more critical review for a certain member who is causing conflicts. Once the TA is involved, everyone must conform to what the majority decides. "

$r2 = $d.Content
$r2.Find.ClearFormatting()
$found2 = $r2.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2.Text = $new2
}
